$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 245, shifting existing rows 245:355 down to 246:356
$ws.Rows("245:245").Insert()

# Populate the newly inserted row 245 with the new record's data
$ws.Range("A245").Value = 2
$ws.Range("B245").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C245").Value = "Coquimbo"
$ws.Range("D245").Value = 44917
$ws.Range("E245").Value = 4
$ws.Range("F245").Value = 100112021
$ws.Range("G245").Value = "Ají"
$ws.Range("H245").Value = "Americana (o)"
$ws.Range("I245").Value = "Primera"
$ws.Range("J245").Value = 200
$ws.Range("K245").Value = 11000
$ws.Range("L245").Value = 12000
$ws.Range("M245").Value = 11500
$ws.Range("N245").Value = "$/caja 25 kilos"
$ws.Range("O245").Value = "Provincia de Limarí"
$ws.Range("P245").Value = 460
$ws.Range("Q245").Value = 25
$ws.Range("R245").Value = "Hortaliza"
